$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1 "Device Types": drop the old "Device Type Names" column, keep just
# the abbreviation + associated tests, rename headers, and tighten the test
# lists (remove the spaces after the commas).
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Clear()

# Write the new data one column to the right (B/C) so that deleting the old
# column A afterwards shifts B->A and C->B, inheriting their best-fit widths.
$ws1.Range("B1").Value = "Device Type"
$ws1.Range("C1").Value = "Tests"

$ws1.Range("B2").Value = "CAM"
$ws1.Range("C2").Value = "Test_1,Test_2,Test_3,Test_4"

$ws1.Range("B3").Value = "INP"
$ws1.Range("C3").Value = "Test_1,Test_2,Test_4"

$ws1.Range("B4").Value = "KEY"
$ws1.Range("C4").Value = "Test_1,Test_3,Test_4"

$ws1.Range("B5").Value = "GSS"
$ws1.Range("C5").Value = "Test_1,Test_2,Test_3"

$ws1.Columns.Item(1).Delete() | Out-Null

$ws1.Range("D5").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet 2 "Test Types": add an "Alt Text" column next to the test names,
# reorder to ascending Test_1..Test_4, and rename the headers.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Clear()

$ws2.Range("A1").Value = "Test"
$ws2.Range("B1").Value = "Alt Text"

$ws2.Range("A2").Value = "Test_1"
$ws2.Range("B2").Value = "Open"

$ws2.Range("A3").Value = "Test_2"
$ws2.Range("B3").Value = "Closed"

$ws2.Range("A4").Value = "Test_3"
$ws2.Range("B4").Value = "On"

$ws2.Range("A5").Value = "Test_4"
$ws2.Range("B5").Value = "Off"

$ws2.Range("F4").Select() | Out-Null

# Device Types becomes the active / selected tab again.
$ws1.Activate() | Out-Null
